$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.875.82"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").Value = "2.589.29"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'519.71"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'139.99"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "2.598.96"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "3.040.67"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "58.830.13"
$ws.Range("D16").Value = "'20.50"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "2.579.66"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'339.20"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "'10.19"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "'6.49"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'66.27"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'0.406"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'7.06"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0724"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'5.96"
$ws.Range("E31").Value = "  -4.45%  "
$ws.Range("D32").Value = "'18.76"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").Value = "'1.57"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'149.02"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'36.25"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.47"
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.832"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").Value = "'0.820"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'275.99"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "'10.74"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'0.588"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "'0.0522"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "1.985.37"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "'0.0220"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "'4.51"
$ws.Range("E51").Value = "  -0.10%  "
